$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: sd/Statement-non-opinion -> sv/Statement-opinion
$ws.Range("I12").Value = "sv"
$ws.Range("J12").Value = "Statement-opinion"

# Row 17: sd/Statement-non-opinion -> aa/Agree/Accept
$ws.Range("I17").Value = "aa"
$ws.Range("J17").Value = "Agree/Accept"

# Row 20: sd/Statement-non-opinion -> aa/Agree/Accept
$ws.Range("I20").Value = "aa"
$ws.Range("J20").Value = "Agree/Accept"

# Row 35: sd/Statement-non-opinion -> sv/Statement-opinion
$ws.Range("I35").Value = "sv"
$ws.Range("J35").Value = "Statement-opinion"

# Row 46: b/Acknowledge (Backchannel) -> %/Uninterpretable
$ws.Range("I46").Value = "%"
$ws.Range("J46").Value = "Uninterpretable"

# Row 55: b/Acknowledge (Backchannel) -> %/Uninterpretable
$ws.Range("I55").Value = "%"
$ws.Range("J55").Value = "Uninterpretable"

# Row 59: b/Acknowledge (Backchannel) -> %/Uninterpretable
$ws.Range("I59").Value = "%"
$ws.Range("J59").Value = "Uninterpretable"

# Row 61: sd/Statement-non-opinion -> aa/Agree/Accept
$ws.Range("I61").Value = "aa"
$ws.Range("J61").Value = "Agree/Accept"

# Row 67: sv/Statement-opinion -> ba/Appreciation
$ws.Range("I67").Value = "ba"
$ws.Range("J67").Value = "Appreciation"

# Row 88: sd/Statement-non-opinion -> sv/Statement-opinion
$ws.Range("I88").Value = "sv"
$ws.Range("J88").Value = "Statement-opinion"

# Row 93: sd/Statement-non-opinion -> sv/Statement-opinion
$ws.Range("I93").Value = "sv"
$ws.Range("J93").Value = "Statement-opinion"

# Row 94: sd/Statement-non-opinion -> aa/Agree/Accept
$ws.Range("I94").Value = "aa"
$ws.Range("J94").Value = "Agree/Accept"

# Row 95: %/Uninterpretable -> sd/Statement-non-opinion
$ws.Range("I95").Value = "sd"
$ws.Range("J95").Value = "Statement-non-opinion"
